# Update countries & provincias Spain
# Refresh the COVID-19 country data table: update the "last updated" timestamp
# and refresh the case counts for the countries whose totals changed. A few
# countries (Chile, Pakistan, Irak, Republica de Chipre) overtook their
# neighbours in the ranking, so their rows now carry different countries'
# data than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 17:22"

function Set-Row {
    param($row, $country, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes)
    $ws.Range("A$row").Value = $country
    $ws.Range("B$row").Value = $total
    $ws.Range("C$row").Value = $nuevos
    $ws.Range("D$row").Value = $activos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $criticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

Set-Row 4  "Estados Unidos" 338999 2326 18002 311310 8702 71 9687
Set-Row 7  "Alemania"       100770 647  28700 70462  3936 24 1608
Set-Row 17 "Austria"        12261  210  3463  8578   250  16 220
Set-Row 19 "Brasil"         11494  240  127   10872  296  9  495

Set-Row 27 "Chile"          4815 344 728  4050 327 3 37
Set-Row 28 "Dinamarca"      4681 312 1378 3116 144 8 187
Set-Row 29 "Chequia"        4591 4   96   4423 84  5 72
Set-Row 30 "India"          4553 264 328  4107 0   0 118

Set-Row 35 "Pakistan"       3658 501 257 3349 17  5 52
Set-Row 36 "Japon"          3654 0   575 2994 69  0 85
Set-Row 37 "Ecuador"        3646 0   100 3366 100 0 180

Set-Row 49 "Grecia"         1755 20 269 1407 90 6 79
Set-Row 52 "Islandia"       1562 76 460 1096 11 2 6

Set-Row 63 "Irak"           1031 70  344 623 0  3 64
Set-Row 64 "Eslovenia"      1021 24  102 889 30 2 30
Set-Row 65 "Moldavia"       965  101 37  909 80 4 19

Set-Row 79 "Republica de Macedonia" 570 15 30 519 15 3 21

Set-Row 85 "Republica de Chipre" 465 19 45 411 11 0 9
Set-Row 86 "Costa Rica"          454 0  16 436 14 0 2
